# EventConfig.xlsx update — "Add Most of Functionality"
#
# Applies the content changes to the Event / EventCondition / EventOption
# sheets:
#   - Event sheet gains a new row (id 14, "获得" / "你获得了一个宝贝~")
#   - Several Event rows get updated Desc text and Effects formulas
#   - The "结婚" row's Effects gains a TAGACTIVE effect
#   - One AppearCondition cell is fixed from AGEBETWEEN,0,0 to AGEBETWEEN,10,20;
#   - Two now-redundant AppearCondition cells are cleared
#   - The Event sheet becomes the active / selected tab (was EventOption)

$wb = $excel.ActiveWorkbook
$wsEvent = $wb.Worksheets.Item("Event")
$wsOption = $wb.Worksheets.Item("EventOption")

# --- Event sheet -----------------------------------------------------

# Options column header switched to reference TbEventOption ids instead of
# the free-form EventCondition list description.
$wsEvent.Range("K2").Value = "(list#sep=;),int#ref=TbEventOption"

# Drop the now-unused AppearCondition values on rows 4 and 5.
$wsEvent.Range("J4").ClearContents()
$wsEvent.Range("J5").ClearContents()

# New row 17: 获得 / 你获得了一个宝贝~ (added first, establishing the new
# shared-string entries in the same order the original author typed them).
$wsEvent.Range("B17").Value = "获得"
$wsEvent.Range("C17").Value = "你获得了一个宝贝~"
$wsEvent.Range("D17").Value = "Normal"
$wsEvent.Range("E17").Value = 14
$wsEvent.Range("F17").Value = 0
$wsEvent.Range("G17").Value = $false
$wsEvent.Range("H17").Value = $false
$wsEvent.Range("I17").Value = $false
$wsEvent.Range("L17").Value = "ADD_RESOURCE,1,0"
$wsEvent.Range("M17").Value = $true

# 佳偶天成 (row 6): better description + fixed age-range appear condition.
$wsEvent.Range("C6").Value = "你找到了很棒的配偶"
$wsEvent.Range("J6").Value = "AGEBETWEEN,10,20;"

# 意外身故 (row 7)
$wsEvent.Range("C7").Value = "你死了"

# 痛失爱子 (row 8): new desc + this event can now recur (IsGenUnique -> TRUE)
$wsEvent.Range("C8").Value = "你儿子死了"
$wsEvent.Range("G8").Value = $true

# 明镜高悬 (row 9): new desc + resource effect now grants the resource.
$wsEvent.Range("C9").Value = "你白了"
$wsEvent.Range("L9").Value = "ADD_RESOURCE,1,0"

# 偶得机缘 (row 10)
$wsEvent.Range("C10").Value = "你得到了机缘"

# 登堂入室 (row 12)
$wsEvent.Range("C12").Value = "你读书很厉害"

# 小有进财 (row 13)
$wsEvent.Range("C13").Value = "你赚了点小钱"

# 富甲一方 (row 14)
$wsEvent.Range("C14").Value = "你赚了很多"

# 仙人指路 (row 16): resource effect now grants the resource.
$wsEvent.Range("L16").Value = "ADD_RESOURCE,1,0"

# 结婚 (row 5): marrying now also activates a tag (typed in last, matching
# the trailing position of this string in the shared-string table).
$wsEvent.Range("L5").Value = "MARRY,0,0;TAGACTIVE,1,0"

# --- Active tab / selection -------------------------------------------
# Event becomes the selected sheet (previously EventOption was selected).
$wsOption.Range("F2").Select()
$wsEvent.Activate()
$wsEvent.Range("L16").Select()
